# Revert "added requirement openpyxl"
# - Insert a new "sex" column (E) into the test_file sheet, pushing
#   customer_type from E to F, and populate the sex values per row.
# - Trim Sheet2 back down to just the first two data rows (rows 4-7 removed).
# - Restore the selection/active-sheet state: Sheet2 becomes the
#   active/selected sheet (cell A3 selected) and test_file loses its
#   former selection (I8) / tabSelected flag.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("test_file")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- test_file: insert "sex" column before "customer_type" -----------------
$ws1.Columns("E:E").Insert()

$ws1.Range("E1").Value = "sex"

$sexValues = @{
    2  = "Male"
    3  = "Other"
    4  = "Female"
    5  = "Male"
    6  = "Female"
    7  = "Male"
    8  = "Female"
    9  = "Male"
    10 = "Female"
    11 = "Male"
    12 = "Male"
    13 = "Male"
    14 = "Female"
    15 = "Female"
    16 = "Female"
    17 = "Female"
    18 = "Female"
    19 = "Male"
    20 = "Female"
    21 = "Male"
    22 = "Male"
    23 = "Female"
    24 = "Female"
    25 = "Female"
    26 = "Male"
    27 = "Male"
    28 = "Male"
    29 = "Male"
    30 = "Female"
}

foreach ($row in 2..30) {
    $ws1.Range("E$row").Value = $sexValues[$row]
}

# --- Sheet2: drop the extra duplicated rows (keep header + 2 data rows) ----
$ws2.Rows("4:7").Delete()

# --- Selection / active sheet state -----------------------------------------
# Reset test_file's lingering selection (was I8) back to the top-left cell,
# then make Sheet2 the active tab with A3 selected.
$ws1.Select()
$ws1.Range("A1").Select()
$ws2.Select()
$ws2.Range("A3").Select()
